# Update the cm020 row: the "Shiny objects" topic for cm020 is replaced by
# "Building Shiny apps (part II)" and the link_it flag is flipped to TRUE.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D21").Value = "Building Shiny apps (part II)"
$ws.Range("C21").Value = $true

# Reflect the new active selection (the author ended up with C22 selected).
$ws.Range("C22").Select()
